$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C1 used a duplicate bold/white font style (no charset) that is identical in
# appearance to the style used by A1/B1. Align it first so that, once the
# "Parent/Child Contract ID" columns are removed, that now-unused duplicate
# font/style entry can be dropped entirely.
$ws.Range("C1").Style = $ws.Range("B1").Style

# Delete columns G:H ("Parent Contract ID" and "Chid Contract ID")
$ws.Range("G1:H1").EntireColumn.Delete()

# Update selection to match target state
$ws.Range("G1").Select()
